# Auto-generated from verified diff analysis (see commit: chore: update Sheets via scheduled runner)
# Refreshes cached marketboard price/profit figures on several 'Profits' worksheets.
# All target cells are plain numeric literals (no formulas in this workbook), so each
# edit is applied as a direct Range.Value write (or ClearContents for cells removed).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 386.3
$ws.Range("I6").Value = 143
$ws.Range("K6").Value = 429
$ws.Range("M6").Value = -317
$ws.Range("H7").Value = 6966.6665
$ws.Range("H10").Value = 25417.666
$ws.Range("I10").Value = 7503
$ws.Range("K10").Value = 7503
$ws.Range("M10").Value = -7210
$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 400
$ws.Range("K12").Value = 400
$ws.Range("M12").Value = -230
$ws.Range("H14").Value = 6966.6665
$ws.Range("H51").Value = 7066.6665
$ws.Range("J51").Value = 7066.6665
$ws.Range("L51").Value = 7066.6665
$ws.Range("N51").Value = -8034.6665
$ws.Range("H86").Value = 3292592.5
$ws.Range("I86").Value = 5149
$ws.Range("J86").Value = 4051233.2
$ws.Range("K86").Value = 5149
$ws.Range("L86").Value = 4051233.2
$ws.Range("M86").Value = -4026
$ws.Range("N86").Value = -4053479.2
$ws.Range("H88").Value = 4099.6
$ws.Range("J88").Value = 2952.8667
$ws.Range("L88").Value = 2952.8667
$ws.Range("N88").Value = -3764.8667
$ws.Range("H89").Value = 3292592.5
$ws.Range("I89").Value = 5149
$ws.Range("J89").Value = 4051233.2
$ws.Range("K89").Value = 25745
$ws.Range("L89").Value = 20256166
$ws.Range("M89").Value = -20129
$ws.Range("N89").Value = -20267398
$ws.Range("H91").Value = 4099.6
$ws.Range("J91").Value = 2952.8667
$ws.Range("L91").Value = 2952.8667
$ws.Range("N91").Value = -5760.8667
$ws.Range("H98").Value = 1054
$ws.Range("I98").Value = 1060
$ws.Range("K98").Value = 1060
$ws.Range("M98").Value = 438
$ws.Range("H100").Value = 4598.811
$ws.Range("I100").Value = 3415.35
$ws.Range("J100").Value = 5991.1177
$ws.Range("K100").Value = 3415.35
$ws.Range("L100").Value = 5991.1177
$ws.Range("M100").Value = -2874.35
$ws.Range("N100").Value = -7073.1177
$ws.Range("H101").Value = 631.1429000000001
$ws.Range("I101").Value = 611.3333
$ws.Range("J101").Value = 750
$ws.Range("K101").Value = 1833.9999
$ws.Range("L101").Value = 2250
$ws.Range("M101").Value = -211.9999
$ws.Range("N101").Value = -5494
$ws.Range("H107").Value = 974.1429000000001
$ws.Range("I107").Value = 803.3333
$ws.Range("K107").Value = 803.3333
$ws.Range("M107").Value = 1116.6667
$ws.Range("H113").Value = 15432.857
$ws.Range("I113").Value = 10825
$ws.Range("K113").Value = 10825
$ws.Range("M113").Value = -7571
$ws.Range("H122").Value = 1054
$ws.Range("I122").Value = 1060
$ws.Range("K122").Value = 3180
$ws.Range("M122").Value = -730
$ws.Range("H131").Value = 10206.652
$ws.Range("I131").Value = 9589.083000000001
$ws.Range("J131").Value = 10880.363
$ws.Range("K131").Value = 28767.249
$ws.Range("L131").Value = 32641.089
$ws.Range("M131").Value = -23727.249
$ws.Range("N131").Value = -42721.089
$ws.Range("H137").Value = 6541.8696
$ws.Range("I137").Value = 6497.8237
$ws.Range("K137").Value = 19493.4711
$ws.Range("M137").Value = -16943.4711
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7719
$ws.Range("J45").Value = 9689.5
$ws.Range("L45").Value = 9689.5
$ws.Range("N45").Value = -10443.5
$ws.Range("H61").Value = 18828.125
$ws.Range("I61").Value = 2593.75
$ws.Range("J61").Value = 35062.5
$ws.Range("K61").Value = 2593.75
$ws.Range("L61").Value = 35062.5
$ws.Range("M61").Value = -2381.75
$ws.Range("N61").Value = -35486.5
$ws.Range("H74").Value = 3828.4092
$ws.Range("I74").Value = 3429.111
$ws.Range("K74").Value = 3429.111
$ws.Range("M74").Value = -2555.111
$ws.Range("H77").Value = 3828.4092
$ws.Range("I77").Value = 3429.111
$ws.Range("K77").Value = 17145.555
$ws.Range("M77").Value = -12777.555
$ws.Range("H88").Value = 5209848
$ws.Range("J88").Value = 8334824.5
$ws.Range("L88").Value = 8334824.5
$ws.Range("N88").Value = -8335636.5
$ws.Range("H91").Value = 5209848
$ws.Range("J91").Value = 8334824.5
$ws.Range("L91").Value = 8334824.5
$ws.Range("N91").Value = -8337632.5
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H130").Value = 208485.8
$ws.Range("J130").Value = 208485.8
$ws.Range("L130").Value = 208485.8
$ws.Range("N130").Value = -218525.8
$ws.Range("H132").Value = 4705.7144
$ws.Range("I132").Value = 2988
$ws.Range("K132").Value = 8964
$ws.Range("M132").Value = -6434
$ws.Range("H136").Value = 18828.125
$ws.Range("I136").Value = 2593.75
$ws.Range("J136").Value = 35062.5
$ws.Range("K136").Value = 7781.25
$ws.Range("L136").Value = 105187.5
$ws.Range("M136").Value = -5231.25
$ws.Range("N136").Value = -110287.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2586811.2
$ws.Range("I94").Value = 2049.7646
$ws.Range("K94").Value = 2049.7646
$ws.Range("M94").Value = -1598.7646
$ws.Range("H132").Value = 52406.777
$ws.Range("J132").Value = 52406.777
$ws.Range("L132").Value = 52406.777
$ws.Range("N132").Value = -62526.777
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2062.3845
$ws.Range("I31").Value = 1900.9166
$ws.Range("K31").Value = 1900.9166
$ws.Range("M31").Value = -1605.9166
$ws.Range("H34").Value = 2062.3845
$ws.Range("I34").Value = 1900.9166
$ws.Range("K34").Value = 1900.9166
$ws.Range("M34").Value = -1698.9166
$ws.Range("H107").Value = 6061.7617
$ws.Range("I107").Value = 649.4286
$ws.Range("J107").Value = 16886.428
$ws.Range("K107").Value = 649.4286
$ws.Range("L107").Value = 16886.428
$ws.Range("M107").Value = 1270.5714
$ws.Range("N107").Value = -20726.428
$ws.Range("H132").Value = 2026.5
$ws.Range("I132").Value = 1971.5294
$ws.Range("K132").Value = 5914.5882
$ws.Range("M132").Value = -3384.5882
$ws.Range("H134").Value = 5473.92
$ws.Range("I134").Value = 5515.1816
$ws.Range("K134").Value = 16545.5448
$ws.Range("M134").Value = -14010.5448
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4089.7036
$ws.Range("I122").Value = 3213.3
$ws.Range("K122").Value = 9639.900000000001
$ws.Range("M122").Value = -7189.900000000001
$ws.Range("H132").Value = 8673.593000000001
$ws.Range("I132").Value = 8673.593000000001
$ws.Range("K132").Value = 26020.779
$ws.Range("M132").Value = -23490.779
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4495.391
$ws.Range("I7").Value = 4341.6665
$ws.Range("J7").Value = 5048.8
$ws.Range("K7").Value = 4341.6665
$ws.Range("L7").Value = 5048.8
$ws.Range("M7").Value = -4229.6665
$ws.Range("N7").Value = -5272.8
$ws.Range("H46").Value = 3127
$ws.Range("I46").Value = 1499.6666
$ws.Range("K46").Value = 1499.6666
$ws.Range("M46").Value = -1311.6666
$ws.Range("H82").Value = 38463170
$ws.Range("I82").Value = 55557264
$ws.Range("K82").Value = 55557264
$ws.Range("M82").Value = -55556903
$ws.Range("H85").Value = 38463170
$ws.Range("I85").Value = 55557264
$ws.Range("K85").Value = 55557264
$ws.Range("M85").Value = -55556016
$ws.Range("H116").Value = 175495.62
$ws.Range("J116").Value = 175495.62
$ws.Range("L116").Value = 175495.62
$ws.Range("N116").Value = -184673.62
$ws.Range("H126").Value = 4495.391
$ws.Range("I126").Value = 4341.6665
$ws.Range("J126").Value = 5048.8
$ws.Range("K126").Value = 13024.9995
$ws.Range("L126").Value = 15146.4
$ws.Range("M126").Value = -10554.9995
$ws.Range("N126").Value = -20086.4
$ws.Range("H132").Value = 3868.8333
$ws.Range("I132").Value = 2654
$ws.Range("K132").Value = 7962
$ws.Range("M132").Value = -5432
$ws.Range("H133").Value = 71367.2
$ws.Range("J133").Value = 71367.2
$ws.Range("L133").Value = 71367.2
$ws.Range("N133").Value = -76427.2
$ws.Range("H136").Value = 5777.6665
$ws.Range("J136").Value = 6999.75
$ws.Range("L136").Value = 20999.25
$ws.Range("N136").Value = -26099.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 35499
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H113").Value = 1405.4375
$ws.Range("I113").Value = 544.2727
$ws.Range("K113").Value = 1632.8181
$ws.Range("M113").Value = 537.1819
$ws.Range("H122").Value = 1663.8
$ws.Range("J122").Value = 1737.25
$ws.Range("L122").Value = 5211.75
$ws.Range("N122").Value = -10111.75
$ws.Range("H140").Value = 76000
$ws.Range("J140").Value = 76000
$ws.Range("L140").Value = 76000
$ws.Range("N140").Value = -86360
